# Auto-generated edit script for Project List.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: row height change (128.95 -> 141.7) ---
$ws.Rows.Item(3).RowHeight = 141.7

# --- Row 7 (height=90.7) ---
# A7: general horiz, no wrap (style like D4)
$ws.Range("A7").Value = "Police Alert"
$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Size = 12
$ws.Range("A7").HorizontalAlignment = 1
$ws.Range("A7").WrapText = $false

# B7: left horiz, no wrap (style like B2/B3)
$ws.Range("B7").Value = "Android"
$ws.Range("B7").Font.Name = "Arial"
$ws.Range("B7").Font.Size = 12
$ws.Range("B7").HorizontalAlignment = -4131
$ws.Range("B7").WrapText = $false

# C7: general horiz, no wrap (style like D4)
$ws.Range("C7").Value = "Umair, Saad Ahmed, Muhammad Ata Jilani"
$ws.Range("C7").Font.Name = "Arial"
$ws.Range("C7").Font.Size = 12
$ws.Range("C7").HorizontalAlignment = 1
$ws.Range("C7").WrapText = $false

# D7: justify horiz, no wrap, dark-gray font color (new style)
$ws.Range("D7").Value = "Police Alert Mobile Application enables citizens of the Karachi Police to receive various Safety, Security and Public services. You can avail various Police services through mobile phone. The efforts to improvise and add new services will continue. "
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("D7").Font.Size = 12
$ws.Range("D7").HorizontalAlignment = -4130
$ws.Range("D7").WrapText = $false
$ws.Range("D7").Font.Color = 3355443

$ws.Rows.Item(7).RowHeight = 90.7

# --- Row 8 (height=15) ---
# A8: general horiz, wrap text (style like A5/D5)
$ws.Range("A8").Value = "Joystick Application"
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Size = 12
$ws.Range("A8").HorizontalAlignment = 1
$ws.Range("A8").WrapText = $true

# B8: left horiz, no wrap (style like B2/B3)
$ws.Range("B8").Value = "Android"
$ws.Range("B8").Font.Name = "Arial"
$ws.Range("B8").Font.Size = 12
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("B8").WrapText = $false

# C8: general horiz, wrap text (style like A5/D5)
$ws.Range("C8").Value = "Umair Shuja, Sarang, Zohaib"
$ws.Range("C8").Font.Name = "Arial"
$ws.Range("C8").Font.Size = 12
$ws.Range("C8").HorizontalAlignment = 1
$ws.Range("C8").WrapText = $true

$ws.Rows.Item(8).RowHeight = 15

# --- Row 9 (height=65.2) ---
# A9: general horiz, wrap text (style like A5/D5)
$ws.Range("A9").Value = "Appoint A Doctor"
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 12
$ws.Range("A9").HorizontalAlignment = 1
$ws.Range("A9").WrapText = $true

# B9: general horiz, wrap text (style like A5/D5)
$ws.Range("B9").Value = "Android studio,java,xml"
$ws.Range("B9").Font.Name = "Arial"
$ws.Range("B9").Font.Size = 12
$ws.Range("B9").HorizontalAlignment = 1
$ws.Range("B9").WrapText = $true

# C9: left horiz, no wrap (style like B2/B3)
$ws.Range("C9").Value = "Atoofa Akber, Hafiz Umer Bin Nasir, Soha Gazdar"
$ws.Range("C9").Font.Name = "Arial"
$ws.Range("C9").Font.Size = 12
$ws.Range("C9").HorizontalAlignment = -4131
$ws.Range("C9").WrapText = $false

# D9: general horiz, wrap text (style like A5/D5)
$ws.Range("D9").Value = "This app alllows user to find Doctors in Karachi and appoint them by Call(if this service provided by them) and also this app allow user to interact with emergency services in case of Emergency."
$ws.Range("D9").Font.Name = "Arial"
$ws.Range("D9").Font.Size = 12
$ws.Range("D9").HorizontalAlignment = 1
$ws.Range("D9").WrapText = $true

$ws.Rows.Item(9).RowHeight = 65.2

# --- Row 10 (height=77.95) ---
# A10: general horiz, wrap text (style like A5/D5)
$ws.Range("A10").Value = "Android Voting System"
$ws.Range("A10").Font.Name = "Arial"
$ws.Range("A10").Font.Size = 12
$ws.Range("A10").HorizontalAlignment = 1
$ws.Range("A10").WrapText = $true

# B10: general horiz, wrap text (style like A5/D5)
$ws.Range("B10").Value = "Android Studio"
$ws.Range("B10").Font.Name = "Arial"
$ws.Range("B10").Font.Size = 12
$ws.Range("B10").HorizontalAlignment = 1
$ws.Range("B10").WrapText = $true

# C10: general horiz, wrap text (style like A5/D5)
$ws.Range("C10").Value = "Junaid Shabbir, Muhammad Ghazali, Hamza Khalid"
$ws.Range("C10").Font.Name = "Arial"
$ws.Range("C10").Font.Size = 12
$ws.Range("C10").HorizontalAlignment = 1
$ws.Range("C10").WrapText = $true

# D10: general horiz, wrap text (style like A5/D5)
$ws.Range("D10").Value = "This application provides is a new technique of casting votes using mobile phones. Android voting system is an application developed for android devices to deploy an easy and flexible way of casting votes anytime and from anywhere."
$ws.Range("D10").Font.Name = "Arial"
$ws.Range("D10").Font.Size = 12
$ws.Range("D10").HorizontalAlignment = 1
$ws.Range("D10").WrapText = $true

$ws.Rows.Item(10).RowHeight = 77.95

# --- Row 11 (height=15) ---
# A11: general horiz, wrap text (style like A5/D5)
$ws.Range("A11").Value = "Shop in budget"
$ws.Range("A11").Font.Name = "Arial"
$ws.Range("A11").Font.Size = 12
$ws.Range("A11").HorizontalAlignment = 1
$ws.Range("A11").WrapText = $true

# B11: left horiz, no wrap (style like B2/B3)
$ws.Range("B11").Value = "Android, JSON, PHP"
$ws.Range("B11").Font.Name = "Arial"
$ws.Range("B11").Font.Size = 12
$ws.Range("B11").HorizontalAlignment = -4131
$ws.Range("B11").WrapText = $false

# C11: general horiz, no wrap (style like D4)
$ws.Range("C11").Value = "Jawahir Qayyum"
$ws.Range("C11").Font.Name = "Arial"
$ws.Range("C11").Font.Size = 12
$ws.Range("C11").HorizontalAlignment = 1
$ws.Range("C11").WrapText = $false

$ws.Rows.Item(11).RowHeight = 15

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 26.666666666666668   # target stored width ~27.5408163265306
$ws.Columns.Item(2).ColumnWidth = 61.83333333333333   # target stored width ~62.6632653061225
$ws.Columns.Item(3).ColumnWidth = 49.0   # target stored width ~49.8724489795918
$ws.Columns.Item(4).ColumnWidth = 42.166666666666664   # target stored width ~43.0612244897959
$ws.Range($ws.Columns.Item(5), $ws.Columns.Item(1025)).ColumnWidth = 9.5   # target stored width ~10.3928571428571

# --- Selection: A1 (was B1) ---
$ws.Range("A1").Select() | Out-Null

